$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 - default font size
$ws.Range("A15").Value = "Font Size Default"

# Row 16 - 14pt font
$ws.Range("A16").Value = "Font Size 14pt"
$ws.Range("A16").Font.Size = 14
$ws.Range("A16").EntireRow.RowHeight = 18.75

# Row 17 - 8pt font
$ws.Range("A17").Value = "Font Size 8pt"
$ws.Range("A17").Font.Size = 8

# Row 18 - 36pt font with wrap text
$ws.Range("A18").Value = "Font Size 36pt"
$ws.Range("A18").Font.Size = 36
$ws.Range("A18").WrapText = $true
$ws.Range("A18").EntireRow.RowHeight = 92.25

# Row 19 - rich text: "Bold" (bold) + " " + "italic" (italic) + " http://www.google.com"
$ws.Range("A19").Value = "Bold italic http://www.google.com"
$ws.Range("A19").Font.Color = 0
$ws.Range("A19").Characters(1, 4).Font.Bold = $true
$ws.Range("A19").Characters(6, 6).Font.Italic = $true

$ws.Range("D20").Select() | Out-Null
